# Generate Report for Handoff
# Updates the localization-status report:
#  - Refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#    timestamps for the rows whose handoff xliff was just (re)generated.
#  - Sets the "Priority" column to "ht" for those same source files.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = 7,8,10,11,12,13

# Overview sheet: "Latest HO Xliff Generate Date" (column G) and the de-de
# sheet's "Latest Handoff Datetime" (column H) shared the same old
# timestamp "2016-08-30 22:22:42"; both move to "2016-08-30 22:22:57".
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-30 22:22:57"
    $wsDeDe.Range("H$r").Value = "2016-08-30 22:22:57"
}

# zh-cn sheet's "Latest Handoff Datetime" (column H) had its own timestamp
# "2016-08-30 22:22:36", which moves to "2016-08-30 22:22:52".
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-30 22:22:52"
}

# zh-cn and de-de "Priority" column (E) gets set to "ht" for these rows.
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
